$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Exam List"
$ws.Name = "Exam List"

# Row 4: change exam name from "Biology" to "XXX" and points from 81 to 100
$ws.Range("B4").Value = "XXX"
$ws.Range("D4").Value = 100

# Row 20: remove the exam name cell entirely and change points from 73 to 100
$ws.Range("B20").ClearContents()
$ws.Range("D20").Value = 100

# Row 35: remove the exam name cell entirely (points stay at 73)
$ws.Range("B35").ClearContents()

# Update the selection / view so no frozen top-left cell and the new
# active selection is B4 (matches the post-edit sheetView state)
$ws.Range("B4").Select()
